$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the "|S*|/n" / k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary labels (A14:A17) and values (B14:B17)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14 (bold, 12pt, vertically centered) then propagate the same
# format to B15:B17 via copy/paste-format so only one new style record
# is minted.
$b14 = $ws.Range("B14")
$fnt = $b14.Font
$fnt.Bold = $true
$b14.VerticalAlignment = -4108
$fnt.Size = 12

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

$ws.Range("A14:B17").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
